$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C143").Value = "mrfso"
$ws.Range("E143").Value = "Can not be produced by LPJ-GUESS"
$ws.Range("F143").Value = "David Warlind"

$ws.Range("C144").Value = "shrubFrac"
$ws.Range("E144").Value = "Can not be produced by LPJ-GUESS"
$ws.Range("F144").Value = "David Warlind"

$ws.Range("C145").Value = "agesno"
$ws.Range("E145").Value = "Can not be produced by LPJ-GUESS"
$ws.Range("F145").Value = "David Warlind"

$ws.Range("C143:C145").WrapText = $true

$ws.Rows.Item(143).RowHeight = 13.8
$ws.Rows.Item(144).RowHeight = 13.8
$ws.Rows.Item(145).RowHeight = 13.8

$ws.Range("F143:F145").Select()
